$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add three new data rows to the bioscreen template table
$ws.Range("A4").Value = "ST22"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 91
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 96

$ws.Range("A5").Value = "other"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 39
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 31

$ws.Range("A6").Value = "other to"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 90
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 10

# Move the active selection, matching the saved view state
$ws.Range("D7").Select()
